$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30, shifting rows 30-38 down to 31-39
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new data
$ws.Cells.Item(30, 1).Value = 11
$ws.Cells.Item(30, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(30, 3).Value = "Bíobío"
$ws.Cells.Item(30, 4).Value = 44876
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(31, 4).NumberFormat
$ws.Cells.Item(30, 5).Value = 8
$ws.Cells.Item(30, 6).Value = 100112022
$ws.Cells.Item(30, 7).Value = "Arveja Verde"
$ws.Cells.Item(30, 8).Value = "Perfection"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 170
$ws.Cells.Item(30, 11).Value = 20000
$ws.Cells.Item(30, 12).Value = 21000
$ws.Cells.Item(30, 13).Value = 20529
$ws.Cells.Item(30, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(30, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(30, 16).Value = 821
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
